$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A datetime corrections (5 rows): shift from 09:59:59 to 11:00:00 ---
$ws.Range("A5").Value = 37347.45833333334
$ws.Range("A17").Value = 37712.45833333334
$ws.Range("A29").Value = 38078.45833333334
$ws.Range("A41").Value = 38443.45833333334
$ws.Range("A59").Value = 38991.45833333334

# --- Refreshed OHLC data: open/high/low/close (C:F) share one value per row ---
$ws.Range("C205:F205").Value = 4373594970000
$ws.Range("C206:F206").Value = 4411934620000
$ws.Range("C210:F210").Value = 4577407590000
$ws.Range("C211:F211").Value = 4566459490000
$ws.Range("C212:F212").Value = 4592275590000
$ws.Range("C213:F213").Value = 4639859400000
$ws.Range("C214:F214").Value = 4641345140000
$ws.Range("C215:F215").Value = 4681223420000
$ws.Range("C216:F216").Value = 4725508480000
$ws.Range("C217:F217").Value = 4680322510000
$ws.Range("C219:F219").Value = 4809150480000
$ws.Range("C221:F221").Value = 4958595660000
$ws.Range("C222:F222").Value = 5004666910000
$ws.Range("C223:F223").Value = 5020790900000
$ws.Range("C224:F224").Value = 5059232680000
$ws.Range("C225:F225").Value = 5094308060000
$ws.Range("C226:F226").Value = 5178041490000
$ws.Range("C227:F227").Value = 5214187690000
$ws.Range("C228:F228").Value = 5235568230000
$ws.Range("C229:F229").Value = 5179738620000
$ws.Range("C230:F230").Value = 5290478980000
$ws.Range("C231:F231").Value = 5390398340000
$ws.Range("C232:F232").Value = 5449356120000
$ws.Range("C233:F233").Value = 5471474170000
$ws.Range("C234:F234").Value = 5507491430000
$ws.Range("C236:F236").Value = 5564521500000
$ws.Range("C239:F239").Value = 5617130550000
$ws.Range("C240:F240").Value = 5647837280000
$ws.Range("C241:F241").Value = 5542014840000
$ws.Range("C242:F242").Value = 5630383690000
$ws.Range("C243:F243").Value = 5704249840000
$ws.Range("C244:F244").Value = 5739159050000
$ws.Range("C245:F245").Value = 5742427260000
$ws.Range("C246:F246").Value = 5825723830000
$ws.Range("C247:F247").Value = 5801917230000
$ws.Range("C248:F248").Value = 5833040250000
$ws.Range("C249:F249").Value = 5855415460000
$ws.Range("C250:F250").Value = 5887405600000
$ws.Range("C251:F251").Value = 5915934540000
$ws.Range("C252:F252").Value = 5940210650000
$ws.Range("C256:F256").Value = 6077524080000
$ws.Range("C258:F258").Value = 6224248910000

# --- New row appended at the end of the sheet ---
$ws.Range("A258").Copy($ws.Range("A259"))  # carry over date/number format & style
$ws.Range("A259").Value = 45078.41666666666
$ws.Range("B259").Value = "ECONOMICS:CZM2"
$ws.Range("C259:F259").Value = 6243183470000
$ws.Range("G259").Value = 0

Write-Output "applied"
